$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 174.66667
$ws.Range("I9").Value = 127.90909
$ws.Range("J9").Value = 248.14285
$ws.Range("K9").Value = 127.90909
$ws.Range("L9").Value = 248.14285
$ws.Range("M9").Value = 41.09090999999999
$ws.Range("N9").Value = -586.14285
$ws.Range("H12").Value = 1201
$ws.Range("I12").Value = 187.14285
$ws.Range("K12").Value = 187.14285
$ws.Range("M12").Value = -17.14285000000001
$ws.Range("H43").Value = 4600
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H58").Value = 357.72726
$ws.Range("I58").Value = 243.5
$ws.Range("K58").Value = 730.5
$ws.Range("M58").Value = -580.5
$ws.Range("H101").Value = 514.75
$ws.Range("J101").Value = 522.1667
$ws.Range("L101").Value = 1566.5001
$ws.Range("N101").Value = -4810.5001
$ws.Range("H112").Value = 1881.5264
$ws.Range("J112").Value = 1880.7778
$ws.Range("L112").Value = 5642.3334
$ws.Range("N112").Value = -7858.3334
$ws.Range("H116").Value = 6557.391
$ws.Range("I116").Value = 4982.5454
$ws.Range("J116").Value = 8001
$ws.Range("K116").Value = 4982.5454
$ws.Range("L116").Value = 8001
$ws.Range("M116").Value = -1540.5454
$ws.Range("N116").Value = -14885
$ws.Range("H135").Value = 7754.8696
$ws.Range("I135").Value = 4542.067
$ws.Range("J135").Value = 13778.875
$ws.Range("K135").Value = 40878.603
$ws.Range("L135").Value = 124009.875
$ws.Range("M135").Value = -38343.603
$ws.Range("N135").Value = -129079.875
$ws.Range("H138").Value = 4656.7036
$ws.Range("I138").Value = 799.75
$ws.Range("J138").Value = 5327.478
$ws.Range("K138").Value = 2399.25
$ws.Range("L138").Value = 15982.434
$ws.Range("M138").Value = 2740.75
$ws.Range("N138").Value = -26262.434
$ws.Range("H141").Value = 7119.0586
$ws.Range("I141").Value = 4573.1787
$ws.Range("J141").Value = 18999.834
$ws.Range("K141").Value = 13719.5361
$ws.Range("L141").Value = 56999.50199999999
$ws.Range("M141").Value = -8539.536100000001
$ws.Range("N141").Value = -67359.50199999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7131.719
$ws.Range("I32").Value = 7883.755
$ws.Range("K32").Value = 7883.755
$ws.Range("M32").Value = -7596.755
$ws.Range("H45").Value = 3712.5
$ws.Range("I45").Value = 2456.25
$ws.Range("J45").Value = 6225
$ws.Range("K45").Value = 2456.25
$ws.Range("L45").Value = 6225
$ws.Range("M45").Value = -2079.25
$ws.Range("N45").Value = -6979
$ws.Range("H61").Value = 3052.2058
$ws.Range("I61").Value = 2540.9644
$ws.Range("J61").Value = 5438
$ws.Range("K61").Value = 2540.9644
$ws.Range("L61").Value = 5438
$ws.Range("M61").Value = -2328.9644
$ws.Range("N61").Value = -5862
$ws.Range("H110").Value = 2032.3125
$ws.Range("I110").Value = 2216.3914
$ws.Range("K110").Value = 2216.3914
$ws.Range("M110").Value = -171.3914
$ws.Range("H136").Value = 3052.2058
$ws.Range("I136").Value = 2540.9644
$ws.Range("J136").Value = 5438
$ws.Range("K136").Value = 7622.8932
$ws.Range("L136").Value = 16314
$ws.Range("M136").Value = -5072.8932
$ws.Range("N136").Value = -21414

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5116.8477
$ws.Range("I134").Value = 3010.9556
$ws.Range("K134").Value = 9032.8668
$ws.Range("M134").Value = -6497.8668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3882.7637
$ws.Range("I31").Value = 4458.3716
$ws.Range("J31").Value = 2875.45
$ws.Range("K31").Value = 4458.3716
$ws.Range("L31").Value = 2875.45
$ws.Range("M31").Value = -4163.3716
$ws.Range("N31").Value = -3465.45
$ws.Range("H34").Value = 3882.7637
$ws.Range("I34").Value = 4458.3716
$ws.Range("J34").Value = 2875.45
$ws.Range("K34").Value = 4458.3716
$ws.Range("L34").Value = 2875.45
$ws.Range("M34").Value = -4256.3716
$ws.Range("N34").Value = -3279.45
$ws.Range("H52").Value = 99782.664
$ws.Range("J52").Value = 99782.664
$ws.Range("L52").Value = 99782.664
$ws.Range("N52").Value = -100370.664
$ws.Range("H107").Value = 4247.8965
$ws.Range("I107").Value = 424.30435
$ws.Range("K107").Value = 424.30435
$ws.Range("M107").Value = 1495.69565
$ws.Range("H132").Value = 774.73334
$ws.Range("I132").Value = 722.9286
$ws.Range("K132").Value = 2168.7858
$ws.Range("M132").Value = 361.2142000000003

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 502.24
$ws.Range("J12").Value = 591.1667
$ws.Range("L12").Value = 1773.5001
$ws.Range("N12").Value = -2119.5001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 1400000
$ws.Range("J21").Value = 800000
$ws.Range("L21").Value = 800000
$ws.Range("N21").Value = -800346
$ws.Range("H24").Value = 3480366.5
$ws.Range("I24").Value = 15000000
$ws.Range("J24").Value = 24476.6
$ws.Range("K24").Value = 15000000
$ws.Range("L24").Value = 24476.6
$ws.Range("M24").Value = -14999827
$ws.Range("N24").Value = -24822.6
$ws.Range("H30").Value = 1400000
$ws.Range("J30").Value = 800000
$ws.Range("L30").Value = 800000
$ws.Range("N30").Value = -800210
$ws.Range("H102").Value = 2837.7896
$ws.Range("I102").Value = 2480.606
$ws.Range("J102").Value = 5195.2
$ws.Range("K102").Value = 2480.606
$ws.Range("L102").Value = 5195.2
$ws.Range("M102").Value = -858.6060000000002
$ws.Range("N102").Value = -8439.200000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000.25
$ws.Range("K22").Value = 1000.25
$ws.Range("M22").Value = -705.25
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000.25
$ws.Range("K27").Value = 1000.25
$ws.Range("M27").Value = -893.25
$ws.Range("H100").Value = 47550.88
$ws.Range("I100").Value = 72157.31
$ws.Range("J100").Value = 3806.111
$ws.Range("K100").Value = 72157.31
$ws.Range("L100").Value = 3806.111
$ws.Range("M100").Value = -71616.31
$ws.Range("N100").Value = -4888.111
$ws.Range("H124").Value = 104138.86
$ws.Range("J124").Value = 104138.86
$ws.Range("L124").Value = 104138.86
$ws.Range("N124").Value = -113958.86
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 5340.6
$ws.Range("I132").Value = 5504.5386
$ws.Range("K132").Value = 16513.6158
$ws.Range("M132").Value = -13983.6158

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11181778
$ws.Range("I81").Value = 25030250
$ws.Range("J81").Value = 103000
$ws.Range("K81").Value = 50060500
$ws.Range("L81").Value = 206000
$ws.Range("M81").Value = -50059439
$ws.Range("N81").Value = -208122
$ws.Range("H84").Value = 11181778
$ws.Range("I84").Value = 25030250
$ws.Range("J84").Value = 103000
$ws.Range("K84").Value = 250302500
$ws.Range("L84").Value = 1030000
$ws.Range("M84").Value = -250297196
$ws.Range("N84").Value = -1040608
$ws.Range("H96").Value = 37305.75
$ws.Range("I96").Value = 55469.4
$ws.Range("K96").Value = 55469.4
$ws.Range("M96").Value = -54096.4
$ws.Range("H100").Value = 751.25
$ws.Range("I100").Value = 765.7857
$ws.Range("K100").Value = 1531.5714
$ws.Range("M100").Value = -990.5714
$ws.Range("H126").Value = 2746
$ws.Range("I126").Value = 2996.3333
$ws.Range("K126").Value = 8988.999899999999
$ws.Range("M126").Value = -6518.999899999999
$ws.Range("H132").Value = 2564.9443
$ws.Range("I132").Value = 2804.6155
$ws.Range("K132").Value = 8413.8465
$ws.Range("M132").Value = -5883.8465
$ws.Range("H136").Value = 3801.3928
$ws.Range("J136").Value = 4067.9546
$ws.Range("L136").Value = 12203.8638
$ws.Range("N136").Value = -17303.8638
